# KIBON-2046: add new fields to institutionen statistik
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "email benachrichtigung kibon" placeholder value (F5).
# (sharedStrings.xml will garbage-collect the old, now-unused string and
# append the corrected one at the end; every later shared-string index
# shifts down by one as a side effect, same as the source diff shows.)
$ws.Range("F5").Value = "{emailBenachrichtigungKiBon}"

# New column AK: "repeatInstitutionenRow" marker, data row only (no header).
# Give it the same border/number-format combo as its neighbour (AJ5) so it
# reuses the existing style slot instead of minting a near-duplicate one.
$ws.Range("AK5").Value = "{repeatInstitutionenRow}"
$ws.Range("AJ5").Copy() | Out-Null
$ws.Range("AK5").PasteSpecial(-4122) | Out-Null

# AE5 picks up a new style: bordered + left-aligned.
$ws.Range("AE5").HorizontalAlignment = -4131

# Column widths / autofit flags for the touched columns.
$ws.Columns("F:F").ColumnWidth = 31.5
$ws.Columns("Y:Y").ColumnWidth = 9.83
$ws.Columns("Z:Z").ColumnWidth = 17
$ws.Columns("AA:AA").ColumnWidth = 17.17
$ws.Columns("AB:AB").ColumnWidth = 14
$ws.Columns("AC:AC").ColumnWidth = 18.83
$ws.Columns("AD:AD").ColumnWidth = 14.67
$ws.Columns("AE:AE").ColumnWidth = 24.83
$ws.Columns("AF:AF").ColumnWidth = 21.33
$ws.Columns("AG:AG").ColumnWidth = 15
$ws.Columns("AK:AK").ColumnWidth = 10.67

# View state: scrolled/selected near the new last column.
$ws.Application.ActiveWindow.ScrollColumn = 21
$ws.Range("AK5").Select()
